# Updated symbol list on Thu Jan 26 04:55:52 UTC 2023 with GitHub Actions
#
# Refreshes the cryptos price/volume sheet. Columns:
#   B=Coin  C=Link  D=Price  E=Volume(1h)  F=Data  G=Hora
#
# Note: D/E hold numeric-looking text ("307.01", "1.86%") stored as plain
# strings (not numbers), exactly like the source data. Assigning a bare
# numeric-looking string via .Value would make Excel coerce it into a real
# number, so those writes are prefixed with a literal leading apostrophe
# (the standard "force text" idiom), which keeps the stored value as exact
# text while still going through the normal Value setter.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - BNB
$ws.Range("D2").Value = "'307.01"
$ws.Range("E2").Value = "'1.86%"

# Row 3 - OKB
$ws.Range("D3").Value = "'36.22"
$ws.Range("E3").Value = "'3.53%"

# Row 4 - HuobiToken
$ws.Range("D4").Value = "'5.101"
$ws.Range("E4").Value = "'1.45%"

# Row 5 - Cronos
$ws.Range("D5").Value = "'0.08122"
$ws.Range("E5").Value = "'2.74%"

# Row 6 - FTXToken
$ws.Range("E6").Value = "'0.39%"

# Row 7 - was KuCoinToken -> GateToken (rows 7-17 shift one position)
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D7").Value = "'4.185"
$ws.Range("E7").Value = "'4.11%"

# Row 8 - was MXToken -> KuCoinToken
$ws.Range("B8").Value = "KuCoinToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D8").Value = "'7.763"
$ws.Range("E8").Value = "'0.20%"

# Row 9 - was LiechtensteinCryptoassetsExchange -> MXToken
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "'0.9303"
$ws.Range("E9").Value = "'0.91%"

# Row 10 - was WazirX -> LiechtensteinCryptoassetsExchange
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.1419"
$ws.Range("E10").Value = "'20.50%"

# Row 11 - was MandalaExchangeToken -> WazirX
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1930"
$ws.Range("E11").Value = "'5.24%"

# Row 12 - was BitrueCoin -> MandalaExchangeToken
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.09264"
$ws.Range("E12").Value = "'0.14%"

# Row 13 - was BitMartToken -> BitrueCoin
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03517"
$ws.Range("E13").Value = "'-0.36%"

# Row 14 - was BitForexToken -> BitMartToken
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09849"
$ws.Range("E14").Value = "'-0.14%"

# Row 15 - was TigerCash -> BitForexToken
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001414"
$ws.Range("E15").Value = "'1.90%"

# Row 16 - was LEO -> TigerCash
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.005747"
$ws.Range("E16").Value = "'-2.63%"

# Row 17 - was GateToken -> LEO
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.605"
$ws.Range("E17").Value = "'2.88%"

# Row 19
$ws.Range("E19").Value = "'-0.12%"

# Row 20
$ws.Range("D20").Value = "'0.1322"
$ws.Range("E20").Value = "'1.10%"

# Row 21
$ws.Range("D21").Value = "'4.888"
$ws.Range("E21").Value = "'-3.12%"

# Row 22
$ws.Range("E22").Value = "'0.41%"

# Row 23
$ws.Range("D23").Value = "'0.04515"
$ws.Range("E23").Value = "'0.38%"

# Row 24
$ws.Range("D24").Value = "'0.001216"
$ws.Range("E24").Value = "'0.14%"

# Row 25
$ws.Range("D25").Value = "'0.004872"
$ws.Range("E25").Value = "'6.65%"

# Row 26
$ws.Range("E26").Value = "'-0.88%"

# Row 39 - One
$ws.Range("D39").Value = "'0.02011"
$ws.Range("E39").Value = "'5.77%"

# Row 40 - IDEX
$ws.Range("D40").Value = "'0.04934"
$ws.Range("E40").Value = "'5.04%"

# Row 41 - was Dexo -> KickToken
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "'0.007649"
$ws.Range("E41").Value = "'1.01%"

# Row 42 - was KickToken -> Dexo
$ws.Range("B42").Value = "Dexo"
$ws.Range("C42").Value = "https://coinranking.com/coin/QkL_pl546+dexo-dexo"
$ws.Range("D42").Value = "'0.01044"
$ws.Range("E42").Value = "'9.25%"

# Row 43 - BKEXToken
$ws.Range("E43").Value = "'4.51%"

# Row 44 - CEJI
$ws.Range("D44").Value = "'0.002100"
$ws.Range("E44").Value = "'-0.56%"

# Row 45 - LocalTraders
$ws.Range("D45").Value = "'0.009997"
$ws.Range("E45").Value = "'-10.43%"

# Row 46 - CoinLion
$ws.Range("D46").Value = "'0.00006437"
$ws.Range("E46").Value = "'7.15%"

# Row 47 - Kangarootoken
$ws.Range("E47").Value = "'0.03%"

# Row 49 - CoinbaseStockToken
$ws.Range("E49").Value = "'-8.73%"

# Row 50 - CryptobidCoin
$ws.Range("E50").Value = "'0.03%"

# Row 51 - SpecialPowerGold
$ws.Range("E51").Value = "'0.03%"
